$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 79, shifting the existing rows 79:105 down to 82:108.
$ws.Rows("79:81").Insert()

# Shared metadata for this block of "Femacal de La Calera" / Membrillo / "$/caja 18 kilos empedrada" rows.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$productoId = 100104
$producto = "Frutos de pepita"
$categoriaId = 100104003
$categoria = "Membrillo"
$variedad = "Champion"
$unidad = "`$/caja 18 kilos empedrada"
$origen = "Región de O'Higgins"
$kgUnidad = 18

function Set-MembrilloRow($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

Set-MembrilloRow 79 45016 "Especial" 70 16000 16000 16000 889
Set-MembrilloRow 80 45016 "Extra (doble especial)" 68 18000 18000 18000 1000
Set-MembrilloRow 81 45016 "Primera" 60 14000 14000 14000 778
